$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3491903333333333
$ws.Range("H2").Value = 1.047571
$ws.Range("I2").Value = 0.008130334326258625
$ws.Range("J2").Value = 0.008130334326258625
$ws.Range("M2").Value = 33.54649666666666
$ws.Range("N2").Value = 100.63949
$ws.Range("O2").Value = 0.1978943147725085
$ws.Range("P2").Value = 0.1978943147725085
$ws.Range("Q2").Value = 11.71411235319889
$ws.Range("R2").Value = 105.42701117879
$ws.Range("S2").Value = 0.001608946940366355
$ws.Range("T2").Value = 0.001608946940366355

# Row 3
$ws.Range("G3").Value = 0.3491903333333333
$ws.Range("H3").Value = 1.047571
$ws.Range("I3").Value = 0.008130334326258625
$ws.Range("J3").Value = 0.008130334326258625
$ws.Range("O3").Value = 0.6322484766686425
$ws.Range("P3").Value = 0.6322484766686425
$ws.Range("Q3").Value = 37.42517666234744
$ws.Range("R3").Value = 336.826589961127
$ws.Range("S3").Value = 0.00514039149258379
$ws.Range("T3").Value = 0.00514039149258379

# Row 4
$ws.Range("G4").Value = 0.3491903333333333
$ws.Range("H4").Value = 1.047571
$ws.Range("I4").Value = 0.008130334326258625
$ws.Range("J4").Value = 0.008130334326258625
$ws.Range("O4").Value = 0.169857208558849
$ws.Range("P4").Value = 0.169857208558849
$ws.Range("Q4").Value = 10.054490081468
$ws.Range("R4").Value = 90.49041073321199
$ws.Range("S4").Value = 0.00138099589330848
$ws.Range("T4").Value = 0.00138099589330848

# Row 5
$ws.Range("I5").Value = 0.801301577139928
$ws.Range("J5").Value = 0.8013015771399279
$ws.Range("M5").Value = 33.54649666666666
$ws.Range("N5").Value = 100.63949
$ws.Range("O5").Value = 0.1978943147725085
$ws.Range("P5").Value = 0.1978943147725085
$ws.Range("Q5").Value = 1154.508083769297
$ws.Range("R5").Value = 10390.57275392367
$ws.Range("S5").Value = 0.1585730265342364
$ws.Range("T5").Value = 0.1585730265342364

# Row 6
$ws.Range("I6").Value = 0.801301577139928
$ws.Range("J6").Value = 0.8013015771399279
$ws.Range("O6").Value = 0.6322484766686425
$ws.Range("P6").Value = 0.6322484766686425
$ws.Range("S6").Value = 0.5066217014989002
$ws.Range("T6").Value = 0.5066217014989001

# Row 7
$ws.Range("I7").Value = 0.801301577139928
$ws.Range("J7").Value = 0.8013015771399279
$ws.Range("O7").Value = 0.169857208558849
$ws.Range("P7").Value = 0.169857208558849
$ws.Range("S7").Value = 0.1361068491067914
$ws.Range("T7").Value = 0.1361068491067914

# Row 8
$ws.Range("I8").Value = 0.1905680885338134
$ws.Range("J8").Value = 0.1905680885338134
$ws.Range("M8").Value = 33.54649666666666
$ws.Range("N8").Value = 100.63949
$ws.Range("O8").Value = 0.1978943147725085
$ws.Range("P8").Value = 0.1978943147725085
$ws.Range("Q8").Value = 274.5687828370899
$ws.Range("R8").Value = 2471.11904553381
$ws.Range("S8").Value = 0.03771234129790574
$ws.Range("T8").Value = 0.03771234129790574

# Row 9
$ws.Range("I9").Value = 0.1905680885338134
$ws.Range("J9").Value = 0.1905680885338134
$ws.Range("O9").Value = 0.6322484766686425
$ws.Range("P9").Value = 0.6322484766686425
$ws.Range("S9").Value = 0.1204863836771585
$ws.Range("T9").Value = 0.1204863836771585

# Row 10
$ws.Range("I10").Value = 0.1905680885338134
$ws.Range("J10").Value = 0.1905680885338134
$ws.Range("O10").Value = 0.169857208558849
$ws.Range("P10").Value = 0.169857208558849
$ws.Range("S10").Value = 0.03236936355874914
$ws.Range("T10").Value = 0.03236936355874914
